# chore: update Sheets via scheduled runner
#
# Refreshes cached market-board figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) on a handful of rows across each crafting-job
# sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Some rows gain a
# previously-absent LeveProfitHQ (column N) value, and a few rows in LTW/WVR
# had their N (and M) figures reset to blank because H-L are now all zero.

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H33").Value = 151.08333
$ws.Range("I33").Value = 151.08333
$ws.Range("K33").Value = 151.08333
$ws.Range("M33").Value = 77.91667000000001
$ws.Range("H40").Value = 2166.6667
$ws.Range("J40").Value = 2500
$ws.Range("L40").Value = 2500
$ws.Range("N40").Value = -2850
$ws.Range("H46").Value = 3262
$ws.Range("I46").Value = 655
$ws.Range("K46").Value = 1965
$ws.Range("M46").Value = -1846
$ws.Range("H60").Value = 3262
$ws.Range("I60").Value = 655
$ws.Range("K60").Value = 1965
$ws.Range("M60").Value = -1481
$ws.Range("H98").Value = 26316492
$ws.Range("I98").Value = 623.86664
$ws.Range("J98").Value = 125001000
$ws.Range("K98").Value = 623.86664
$ws.Range("L98").Value = 125001000
$ws.Range("M98").Value = 874.13336
$ws.Range("N98").Value = -125003996
$ws.Range("H113").Value = 2856
$ws.Range("I113").Value = 2763
$ws.Range("J113").Value = 3004.8
$ws.Range("K113").Value = 2763
$ws.Range("L113").Value = 3004.8
$ws.Range("M113").Value = 491
$ws.Range("N113").Value = -9512.799999999999
$ws.Range("H116").Value = 18968.5
$ws.Range("I116").Value = 22361
$ws.Range("J116").Value = 2006
$ws.Range("K116").Value = 22361
$ws.Range("L116").Value = 2006
$ws.Range("M116").Value = -18919
$ws.Range("N116").Value = -8890
$ws.Range("H122").Value = 26316492
$ws.Range("I122").Value = 623.86664
$ws.Range("J122").Value = 125001000
$ws.Range("K122").Value = 1871.59992
$ws.Range("L122").Value = 375003000
$ws.Range("M122").Value = 578.4000800000001
$ws.Range("N122").Value = -375007900
$ws.Range("H138").Value = 2810.1235
$ws.Range("I138").Value = 1470.0571
$ws.Range("J138").Value = 3829.739
$ws.Range("K138").Value = 4410.1713
$ws.Range("L138").Value = 11489.217
$ws.Range("M138").Value = 729.8287
$ws.Range("N138").Value = -21769.217

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 562269.0600000001
$ws.Range("I32").Value = 610864.3
$ws.Range("J32").Value = 27721.334
$ws.Range("K32").Value = 610864.3
$ws.Range("L32").Value = 27721.334
$ws.Range("M32").Value = -610577.3
$ws.Range("N32").Value = -28295.334
$ws.Range("H61").Value = 2779
$ws.Range("I61").Value = 2212.6296
$ws.Range("K61").Value = 2212.6296
$ws.Range("M61").Value = -2000.6296
$ws.Range("H74").Value = 3368.3
$ws.Range("I74").Value = 2722.875
$ws.Range("J74").Value = 5950
$ws.Range("K74").Value = 2722.875
$ws.Range("L74").Value = 5950
$ws.Range("M74").Value = -1848.875
$ws.Range("N74").Value = -7698
$ws.Range("H77").Value = 3368.3
$ws.Range("I77").Value = 2722.875
$ws.Range("J77").Value = 5950
$ws.Range("K77").Value = 13614.375
$ws.Range("L77").Value = 29750
$ws.Range("M77").Value = -9246.375
$ws.Range("N77").Value = -38486
$ws.Range("H102").Value = 4546.6665
$ws.Range("I102").Value = 5236.6665
$ws.Range("K102").Value = 5236.6665
$ws.Range("M102").Value = -3614.6665
$ws.Range("H136").Value = 2779
$ws.Range("I136").Value = 2212.6296
$ws.Range("K136").Value = 6637.888800000001
$ws.Range("M136").Value = -4087.888800000001

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H99").Value = 1495.3334
$ws.Range("I99").Value = 1289.8182
$ws.Range("J99").Value = 2060.5
$ws.Range("K99").Value = 1289.8182
$ws.Range("L99").Value = 2060.5
$ws.Range("M99").Value = 208.1818000000001
$ws.Range("N99").Value = -5056.5
$ws.Range("H106").Value = 94750
$ws.Range("J106").Value = 94750
$ws.Range("L106").Value = 94750
$ws.Range("N106").Value = -97274
$ws.Range("H134").Value = 1839.7462
$ws.Range("I134").Value = 1368.5116
$ws.Range("J134").Value = 2684.0417
$ws.Range("K134").Value = 4105.5348
$ws.Range("L134").Value = 8052.125100000001
$ws.Range("M134").Value = -1570.5348
$ws.Range("N134").Value = -13122.1251

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H29").Value = 5021
$ws.Range("J29").Value = 5021
$ws.Range("L29").Value = 5021
$ws.Range("N29").Value = -5607
$ws.Range("H58").Value = 1079.2609
$ws.Range("I58").Value = 747.53125
$ws.Range("J58").Value = 1837.5
$ws.Range("K58").Value = 747.53125
$ws.Range("L58").Value = 1837.5
$ws.Range("M58").Value = -544.53125
$ws.Range("N58").Value = -2243.5
$ws.Range("H132").Value = 19775178
$ws.Range("I132").Value = 28572522
$ws.Range("J132").Value = 6945717.5
$ws.Range("K132").Value = 85717566
$ws.Range("L132").Value = 20837152.5
$ws.Range("M132").Value = -85715036
$ws.Range("N132").Value = -20842212.5
$ws.Range("H134").Value = 3393.9424
$ws.Range("I134").Value = 3317.4883
$ws.Range("J134").Value = 3759.2222
$ws.Range("K134").Value = 9952.464899999999
$ws.Range("L134").Value = 11277.6666
$ws.Range("M134").Value = -7417.464899999999
$ws.Range("N134").Value = -16347.6666
$ws.Range("H136").Value = 1079.2609
$ws.Range("I136").Value = 747.53125
$ws.Range("J136").Value = 1837.5
$ws.Range("K136").Value = 2242.59375
$ws.Range("L136").Value = 5512.5
$ws.Range("M136").Value = 307.40625
$ws.Range("N136").Value = -10612.5

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H5").Value = 1965.2
$ws.Range("I5").Value = 614.8570999999999
$ws.Range("J5").Value = 2692.3076
$ws.Range("K5").Value = 1844.5713
$ws.Range("L5").Value = 8076.9228
$ws.Range("M5").Value = -1732.5713
$ws.Range("N5").Value = -8300.9228
$ws.Range("H26").Value = 8079.731
$ws.Range("I26").Value = 74.59999999999999
$ws.Range("J26").Value = 9985.714
$ws.Range("K26").Value = 223.8
$ws.Range("L26").Value = 29957.142
$ws.Range("M26").Value = 64.20000000000002
$ws.Range("N26").Value = -30533.142
$ws.Range("H44").Value = 281.625
$ws.Range("I44").Value = 281.625
$ws.Range("K44").Value = 844.875
$ws.Range("M44").Value = -446.875
$ws.Range("H80").Value = 3741.923
$ws.Range("J80").Value = 4063.9
$ws.Range("L80").Value = 12191.7
$ws.Range("N80").Value = -14063.7
$ws.Range("H83").Value = 3741.923
$ws.Range("J83").Value = 4063.9
$ws.Range("L83").Value = 36575.1
$ws.Range("N83").Value = -45935.1
$ws.Range("H109").Value = 3569.9565
$ws.Range("I109").Value = 644.9
$ws.Range("J109").Value = 5820
$ws.Range("K109").Value = 1934.7
$ws.Range("L109").Value = 17460
$ws.Range("M109").Value = -894.6999999999998
$ws.Range("N109").Value = -19540
$ws.Range("H132").Value = 2443.7896
$ws.Range("I132").Value = 2112.8333
$ws.Range("J132").Value = 2596.5386
$ws.Range("K132").Value = 19015.4997
$ws.Range("L132").Value = 23368.8474
$ws.Range("M132").Value = -16485.4997
$ws.Range("N132").Value = -28428.8474
$ws.Range("H135").Value = 1965.2
$ws.Range("I135").Value = 614.8570999999999
$ws.Range("J135").Value = 2692.3076
$ws.Range("K135").Value = 5533.7139
$ws.Range("L135").Value = 24230.7684
$ws.Range("M135").Value = -2998.7139
$ws.Range("N135").Value = -29300.7684

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H33").Value = 1750
$ws.Range("I33").Value = 500
$ws.Range("J33").Value = 3000
$ws.Range("K33").Value = 500
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = -248
$ws.Range("N33").Value = -3504
$ws.Range("H132").Value = 3392.9524
$ws.Range("I132").Value = 3308.4285
$ws.Range("J132").Value = 3562
$ws.Range("K132").Value = 9925.2855
$ws.Range("L132").Value = 10686
$ws.Range("M132").Value = -7395.2855
$ws.Range("N132").Value = -15746

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 12710.7
$ws.Range("I22").Value = 1700
$ws.Range("J22").Value = 17429.572
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 17429.572
$ws.Range("M22").Value = -1405
$ws.Range("N22").Value = -18019.572
$ws.Range("H26").Value = 50006.668
$ws.Range("I26").Value = 10000
$ws.Range("J26").Value = 70010
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 70010
$ws.Range("M26").Value = -9705
$ws.Range("N26").Value = -70600
$ws.Range("H27").Value = 12710.7
$ws.Range("I27").Value = 1700
$ws.Range("J27").Value = 17429.572
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 17429.572
$ws.Range("M27").Value = -1593
$ws.Range("N27").Value = -17643.572
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value = 27375
$ws.Range("I93").Value = 51000
$ws.Range("J93").Value = 3750
$ws.Range("K93").Value = 51000
$ws.Range("L93").Value = 3750
$ws.Range("M93").Value = -49752
$ws.Range("N93").Value = -6246
$ws.Range("H132").Value = 2481.9285
$ws.Range("I132").Value = 2348.8
$ws.Range("J132").Value = 2677.7058
$ws.Range("K132").Value = 7046.400000000001
$ws.Range("L132").Value = 8033.117400000001
$ws.Range("M132").Value = -4516.400000000001
$ws.Range("N132").Value = -13093.1174

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H123").Value = 31204.143
$ws.Range("J123").Value = 98429
$ws.Range("L123").Value = 98429
$ws.Range("N123").Value = -108229
$ws.Range("H132").Value = 2826656.8
$ws.Range("I132").Value = 1851.4445
$ws.Range("J132").Value = 7248091.5
$ws.Range("K132").Value = 5554.333500000001
$ws.Range("L132").Value = 21744274.5
$ws.Range("M132").Value = -3024.333500000001
$ws.Range("N132").Value = -21749334.5
